$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text for the MODEL_CONDITION column (currently in E1, before the column shift)
$ws.Range("E1").Value = "MODELCONDITION"

# Delete the entire first column (A), shifting B:F left to A:E
$ws.Range("A1").EntireColumn.Delete()
